# The class diagram's "User" class listed its "role" attribute as an
# Int; change that attribute's type to Boolean.
#
# "Int: role" lives in a single run. We only want to change the "Int"
# part to "Boolean" while leaving ": role" untouched, which (per the
# target diff) ends up as two separate runs sharing identical
# formatting: "Boolean" and ": role".

$d = $word.ActiveDocument

$match = $d.Content
$found = $match.Find.Execute("Int: role", $true, $false, $false, $false,
                              $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $match.Start

    # Range covering just the "Int" portion of the match.
    $typeRange = $d.Range($start, $start + 3)
    $typeRange.Text = "Boolean"

    # Range covering the just-inserted "Boolean" text.
    $newTypeRange = $d.Range($start, $start + 7)

    # Toggling (and restoring) a character-formatting property forces
    # the engine to keep "Boolean" as its own run instead of folding it
    # back into a single run together with the following ": role" text
    # -- matching the two-run structure produced by Word in this case.
    $newTypeRange.Bold = 1
    $newTypeRange.Bold = 0
}
